$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.806.84"
$ws.Range("E2").Value = "  -0.78%  "

$ws.Range("D3").Value = "'2.922.30"
$ws.Range("E3").Value = "  +0.95%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'354.69"
$ws.Range("E5").Value = "  +0.55%  "

$ws.Range("D6").Value = "'109.92"
$ws.Range("E6").Value = "  -1.98%  "

$ws.Range("D7").Value = "'0.571"
$ws.Range("E7").Value = "  +2.04%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").Value = "'0.630"
$ws.Range("E9").Value = "  +0.76%  "

$ws.Range("D10").Value = "'39.28"
$ws.Range("E10").Value = "  -2.18%  "

$ws.Range("D11").Value = "'0.0884"
$ws.Range("E11").Value = "  +3.03%  "

$ws.Range("D12").Value = "'0.137"
$ws.Range("E12").Value = "  +0.52%  "

$ws.Range("D13").Value = "'19.68"
$ws.Range("E13").Value = "  -1.70%  "

$ws.Range("D14").Value = "'7.88"
$ws.Range("E14").Value = "  +0.97%  "

$ws.Range("D15").Value = "'3.377.33"
$ws.Range("E15").Value = "  +0.91%  "

$ws.Range("D16").Value = "'2.903.38"
$ws.Range("E16").Value = "  -0.73%  "

$ws.Range("D17").Value = "'0.976"
$ws.Range("E17").Value = "  -1.60%  "

$ws.Range("D18").Value = "'51.770.85"
$ws.Range("E18").Value = "  -0.85%  "

$ws.Range("D19").Value = "'3.27"
$ws.Range("E19").Value = "  -2.40%  "

$ws.Range("D20").Value = "'7.52"
$ws.Range("E20").Value = "  -2.75%  "

$ws.Range("D21").Value = "'13.92"
$ws.Range("E21").Value = "  -3.85%  "

$ws.Range("D22").Value = "'0.0₃0981"
$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").Value = "'70.71"
$ws.Range("E23").Value = "  -0.53%  "

$ws.Range("D24").Value = "'270.23"
$ws.Range("E24").Value = "  -0.09%  "

$ws.Range("D25").Value = "'2.81"
$ws.Range("E25").Value = "  +0.54%  "

$ws.Range("D26").Value = "'0.181"
$ws.Range("E26").Value = "  +10.73%  "

$ws.Range("D27").Value = "'27.07"
$ws.Range("E27").Value = "  +2.12%  "

$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("D29").Value = "'7.33"
$ws.Range("E29").Value = "  +12.99%  "

$ws.Range("E30").Value = "  +11.30%  "

$ws.Range("D31").Value = "'10.53"
$ws.Range("E31").Value = "  +0.23%  "

$ws.Range("E32").Value = "  -0.72%  "

$ws.Range("D33").Value = "'6.01"
$ws.Range("E33").Value = "  -2.30%  "

$ws.Range("D34").Value = "'52.22"
$ws.Range("E34").Value = "  -2.05%  "

$ws.Range("D35").Value = "'0.0439"
$ws.Range("E35").Value = "  -4.89%  "

$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("D37").Value = "'1.92"
$ws.Range("E37").Value = "  -14.85%  "

$ws.Range("D38").Value = "'3.22"
$ws.Range("E38").Value = "  -2.95%  "

$ws.Range("D39").Value = "'18.38"
$ws.Range("E39").Value = "  -1.75%  "

$ws.Range("E40").Value = "  -2.08%  "

$ws.Range("D41").Value = "'2.75"
$ws.Range("E41").Value = "  +4.41%  "

$ws.Range("E42").Value = "  +2.25%  "

$ws.Range("D43").Value = "'22.87"
$ws.Range("E43").Value = "  +0.40%  "

$ws.Range("D44").Value = "'121.37"
$ws.Range("E44").Value = "  -0.40%  "

$ws.Range("E45").Value = "  -1.46%  "

$ws.Range("E46").Value = "  +0.99%  "

$ws.Range("D47").Value = "'3.45"
$ws.Range("E47").Value = "  -4.29%  "

$ws.Range("D48").Value = "'2.140.14"
$ws.Range("E48").Value = "  -3.10%  "

$ws.Range("E49").Value = "  -7.28%  "

$ws.Range("D50").Value = "'0.0330"
$ws.Range("E50").Value = "  +2.82%  "

$ws.Range("E51").Value = "  -1.86%  "
